# StrategicDataManagementOrderData.xlsx
# The shipment carrier previously recorded as "DHS" was actually "DHL" —
# correct every tracking number in column J (Tracking number) that was
# mis-entered with the "DHS" prefix so it reads "DHL" instead, across all
# data rows on Sheet1. This is a dataset cleanup for data profiling and
# preprocessing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J holds "Tracking number" (header in J1); data starts at row 2.
# Walk every used row and fix any tracking number that still begins with
# the incorrect "DHS" carrier code.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "DHS*") {
        $newVal = "DHL" + $val.Substring(3)
        $cell.Value = $newVal
    }
}

# Reflect the reviewer's last on-screen selection when the file was saved.
$ws.Range("K34").Select()
